$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "247.47"
Set-TextValue $ws.Range("D3") "21.76"
Set-TextValue $ws.Range("D4") "5.482"
Set-TextValue $ws.Range("D5") "0.05660"
Set-TextValue $ws.Range("D7") "6.434"
Set-TextValue $ws.Range("D8") "0.8001"
Set-TextValue $ws.Range("D10") "0.1431"
Set-TextValue $ws.Range("D12") "0.03160"
Set-TextValue $ws.Range("D13") "0.02954"
Set-TextValue $ws.Range("D14") "0.09280"
Set-TextValue $ws.Range("D15") "0.001643"
Set-TextValue $ws.Range("D16") "3.219"
Set-TextValue $ws.Range("D17") "0.04731"
Set-TextValue $ws.Range("B18") "TigerCash"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D18") "0.006404"
Set-TextValue $ws.Range("E18") "17TigerCashTCH"
Set-TextValue $ws.Range("B19") "HotbitToken"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D19") "0.005011"
Set-TextValue $ws.Range("E19") "18HotbitTokenHTB"
Set-TextValue $ws.Range("B20") "BitKan"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D20") "0.001046"
Set-TextValue $ws.Range("E20") "19BitKanKAN"
Set-TextValue $ws.Range("B21") "NitroEx"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D21") "0.0001502"
Set-TextValue $ws.Range("E21") "20NitroExNTX"
Set-TextValue $ws.Range("B22") "UpBots"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue $ws.Range("D22") "0.0003204"
Set-TextValue $ws.Range("E22") "21UpBotsUBXT"
Set-TextValue $ws.Range("B23") "LEO"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D23") "3.917"
Set-TextValue $ws.Range("E23") "22LEOLEO"
Set-TextValue $ws.Range("B24") "BTSEToken"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D24") "2.090"
Set-TextValue $ws.Range("E24") "23BTSETokenBTSE"
Set-TextValue $ws.Range("B25") "One"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D25") "0.0005857"
Set-TextValue $ws.Range("E25") "24OneONE"
Set-TextValue $ws.Range("D26") "0.3271"
Set-TextValue $ws.Range("D40") "0.04085"
Set-TextValue $ws.Range("D41") "0.006900"
Set-TextValue $ws.Range("E41") "40KickTokenKICKBestin24h"
Set-TextValue $ws.Range("D42") "0.1039"
Set-TextValue $ws.Range("D43") "0.002974"
Set-TextValue $ws.Range("D44") "0.009124"
Set-TextValue $ws.Range("D45") "0.00005816"
Set-TextValue $ws.Range("D48") "0.009148"
